# fix(publipostage): Correct status name
#
# Updates the wording of a few status-related labels used throughout the
# "Liste essais cliniques identifiés" sheet:
#   - "bleu" -> "noir"
#   - "pas de résultat ni de publication" -> "pas de résultat postés ni publiés"
#   - "résultat et / ou publication posté" -> "résultat postés ou publiés"
#   - "résultat et / ou publication posté dans les 12 mois"
#       -> "résultat postés ou publiés dans les 12 mois"

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $cells = $ws.Cells

    # Order matters: replace the longer/more specific strings first so that
    # the shorter substring replacement does not clobber them.
    $cells.Replace(
        "résultat et / ou publication posté dans les 12 mois",
        "résultat postés ou publiés dans les 12 mois",
        [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole
    ) | Out-Null

    $cells.Replace(
        "résultat et / ou publication posté",
        "résultat postés ou publiés",
        [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole
    ) | Out-Null

    $cells.Replace(
        "pas de résultat ni de publication",
        "pas de résultat postés ni publiés",
        [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole
    ) | Out-Null

    $cells.Replace(
        "bleu",
        "noir",
        [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole
    ) | Out-Null
}
